$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 16.71895933333333
$ws.Range("H2").Value = 50.156878
$ws.Range("I2").Value = 0.02912144738161902
$ws.Range("J2").Value = 0.03059269312988411
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 2038.521628905614
$ws.Range("R2").Value = 18346.69466015052
$ws.Range("S2").Value = 0.006646190351280707
$ws.Range("T2").Value = 0.007402576461546661

# Row 3
$ws.Range("G3").Value = 16.71895933333333
$ws.Range("H3").Value = 50.156878
$ws.Range("I3").Value = 0.02912144738161902
$ws.Range("J3").Value = 0.03059269312988411
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 2472.973668087247
$ws.Range("R3").Value = 22256.76301278522
$ws.Range("S3").Value = 0.00806263396902801
$ws.Range("T3").Value = 0.00898022194409347

# Row 4
$ws.Range("G4").Value = 16.71895933333333
$ws.Range("H4").Value = 50.156878
$ws.Range("I4").Value = 0.02912144738161902
$ws.Range("J4").Value = 0.03059269312988411
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 1396.116186415247
$ws.Range("R4").Value = 12565.04567773722
$ws.Range("S4").Value = 0.004551756427721207
$ws.Range("T4").Value = 0.005069780311671305

# Row 5
$ws.Range("G5").Value = 16.71895933333333
$ws.Range("H5").Value = 50.156878
$ws.Range("I5").Value = 0.02912144738161902
$ws.Range("J5").Value = 0.03059269312988411
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 1522.569921086692
$ws.Range("R5").Value = 9135.419526520151
$ws.Range("S5").Value = 0.004964033432458191
$ws.Range("T5").Value = 0.003685985012889346

# Row 6
$ws.Range("G6").Value = 16.71895933333333
$ws.Range("H6").Value = 50.156878
$ws.Range("I6").Value = 0.02912144738161902
$ws.Range("J6").Value = 0.03059269312988411
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 1501.95824465438
$ws.Range("R6").Value = 13517.62420188942
$ws.Range("S6").Value = 0.004896833201130896
$ws.Range("T6").Value = 0.005454129399683328

# Row 7
$ws.Range("I7").Value = 0.2708539632042961
$ws.Range("J7").Value = 0.2845377865576845
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 18959.96634477947
$ws.Range("R7").Value = 170639.6971030152
$ws.Range("S7").Value = 0.06181516231884668
$ws.Range("T7").Value = 0.06885018956160407

# Row 8
$ws.Range("I8").Value = 0.2708539632042961
$ws.Range("J8").Value = 0.2845377865576845
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("S8").Value = 0.07498927974834096
$ws.Range("T8").Value = 0.08352361996770638

# Row 9
$ws.Range("I9").Value = 0.2708539632042961
$ws.Range("J9").Value = 0.2845377865576845
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 12985.05521476641
$ws.Range("R9").Value = 116865.4969328976
$ws.Range("S9").Value = 0.04233516458962407
$ws.Range("T9").Value = 0.04715322257155415

# Row 10
$ws.Range("I10").Value = 0.2708539632042961
$ws.Range("J10").Value = 0.2845377865576845
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 14161.18134438192
$ws.Range("R10").Value = 84967.08806629149
$ws.Range("S10").Value = 0.04616968762028534
$ws.Range("T10").Value = 0.03428276197847462

# Row 11
$ws.Range("I11").Value = 0.2708539632042961
$ws.Range("J11").Value = 0.2845377865576845
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 13969.47541106007
$ws.Range("R11").Value = 125725.2786995406
$ws.Range("S11").Value = 0.04554466892719896
$ws.Range("T11").Value = 0.05072799247834529

# Row 12
$ws.Range("G12").Value = 194.8548433333333
$ws.Range("H12").Value = 584.56453
$ws.Range("I12").Value = 0.3394024086099587
$ws.Range("J12").Value = 0.3565493705749576
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 23758.40533567589
$ws.Range("R12").Value = 213825.648021083
$ws.Range("S12").Value = 0.07745950892292262
$ws.Range("T12").Value = 0.08627497967543131

# Row 13
$ws.Range("G13").Value = 194.8548433333333
$ws.Range("H13").Value = 584.56453
$ws.Range("I13").Value = 0.3394024086099587
$ws.Range("J13").Value = 0.3565493705749576
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 28821.82359890496
$ws.Range("R13").Value = 259396.4123901447
$ws.Range("S13").Value = 0.0939677672256015
$ws.Range("T13").Value = 0.1046620010927452

# Row 14
$ws.Range("G14").Value = 194.8548433333333
$ws.Range("H14").Value = 584.56453
$ws.Range("I14").Value = 0.3394024086099587
$ws.Range("J14").Value = 0.3565493705749576
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 16271.34771700147
$ws.Range("R14").Value = 146442.1294530132
$ws.Range("S14").Value = 0.05304946126920671
$ws.Range("T14").Value = 0.05908688625108185

# Row 15
$ws.Range("G15").Value = 194.8548433333333
$ws.Range("H15").Value = 584.56453
$ws.Range("I15").Value = 0.3394024086099587
$ws.Range("J15").Value = 0.3565493705749576
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 17745.13099304504
$ws.Range("R15").Value = 106470.7859582703
$ws.Range("S15").Value = 0.05785443564388536
$ws.Range("T15").Value = 0.04295913506910667

# Row 16
$ws.Range("G16").Value = 194.8548433333333
$ws.Range("H16").Value = 584.56453
$ws.Range("I16").Value = 0.3394024086099587
$ws.Range("J16").Value = 0.3565493705749576
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 17504.9076094013
$ws.Range("R16").Value = 157544.1684846117
$ws.Range("S16").Value = 0.0570712355483425
$ws.Range("T16").Value = 0.06356636848659256

# Row 17
$ws.Range("G17").Value = 82.82950199999999
$ws.Range("H17").Value = 165.659004
$ws.Range("I17").Value = 0.1442742299952585
$ws.Range("J17").Value = 0.1010420758958371
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 10099.29673085798
$ws.Range("R17").Value = 60595.78038514785
$ws.Range("S17").Value = 0.0329267286329376
$ws.Range("T17").Value = 0.02444935754680872

# Row 18
$ws.Range("G18").Value = 82.82950199999999
$ws.Range("H18").Value = 165.659004
$ws.Range("I18").Value = 0.1442742299952585
$ws.Range("J18").Value = 0.1010420758958371
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 12251.67029256366
$ws.Range("R18").Value = 73510.02175538195
$ws.Range("S18").Value = 0.03994411034491861
$ws.Range("T18").Value = 0.02966003232811795

# Row 19
$ws.Range("G19").Value = 82.82950199999999
$ws.Range("H19").Value = 165.659004
$ws.Range("I19").Value = 0.1442742299952585
$ws.Range("J19").Value = 0.1010420758958371
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 6916.675024405271
$ws.Range("R19").Value = 41500.05014643163
$ws.Range("S19").Value = 0.02255042976160397
$ws.Range("T19").Value = 0.01674455808294683

# Row 20
$ws.Range("G20").Value = 82.82950199999999
$ws.Range("H20").Value = 165.659004
$ws.Range("I20").Value = 0.1442742299952585
$ws.Range("J20").Value = 0.1010420758958371
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 7543.155396780675
$ws.Range("R20").Value = 30172.6215871227
$ws.Range("S20").Value = 0.0245929431924688
$ws.Range("T20").Value = 0.01217413504074508

# Row 21
$ws.Range("G21").Value = 82.82950199999999
$ws.Range("H21").Value = 165.659004
$ws.Range("I21").Value = 0.1442742299952585
$ws.Range("J21").Value = 0.1010420758958371
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 7441.040494756258
$ws.Range("R21").Value = 44646.24296853755
$ws.Range("S21").Value = 0.0242600180633295
$ws.Range("T21").Value = 0.01801399289721856

# Row 22
$ws.Range("G22").Value = 124.2078576666667
$ws.Range("H22").Value = 372.623573
$ws.Range("I22").Value = 0.2163479508088675
$ws.Range("J22").Value = 0.2272780738416368
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 15144.50745918815
$ws.Range("R22").Value = 136300.5671326933
$ws.Range("S22").Value = 0.04937562492490744
$ws.Range("T22").Value = 0.05499493988655382

# Row 23
$ws.Range("G23").Value = 124.2078576666667
$ws.Range("H23").Value = 372.623573
$ws.Range("I23").Value = 0.2163479508088675
$ws.Range("J23").Value = 0.2272780738416368
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 18372.12204750036
$ws.Range("R23").Value = 165349.0984275033
$ws.Range("S23").Value = 0.05989861405110557
$ws.Range("T23").Value = 0.06671552378401852

# Row 24
$ws.Range("G24").Value = 124.2078576666667
$ws.Range("H24").Value = 372.623573
$ws.Range("I24").Value = 0.2163479508088675
$ws.Range("J24").Value = 0.2272780738416368
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 10371.97334541403
$ws.Range("R24").Value = 93347.76010872627
$ws.Range("S24").Value = 0.03381573596991408
$ws.Range("T24").Value = 0.03766421933318926

# Row 25
$ws.Range("G25").Value = 124.2078576666667
$ws.Range("H25").Value = 372.623573
$ws.Range("I25").Value = 0.2163479508088675
$ws.Range("J25").Value = 0.2272780738416368
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 11311.41862812217
$ws.Range("R25").Value = 67868.51176873302
$ws.Range("S25").Value = 0.03687860863457301
$ws.Range("T25").Value = 0.02738377985821364

# Row 26
$ws.Range("G26").Value = 124.2078576666667
$ws.Range("H26").Value = 372.623573
$ws.Range("I26").Value = 0.2163479508088675
$ws.Range("J26").Value = 0.2272780738416368
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 11158.29114443533
$ws.Range("R26").Value = 100424.620299918
$ws.Range("S26").Value = 0.0363793672283674
$ws.Range("T26").Value = 0.04051961097966159
